$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.2503572190582515
$ws.Range("J2").Value = 0.2503572190582515
$ws.Range("M2").Value = 0.5134303333333333
$ws.Range("N2").Value = 1.540291
$ws.Range("O2").Value = 0.03326489761800302
$ws.Range("P2").Value = 0.03326489761800301
$ws.Range("Q2").Value = 1.037622851027111
$ws.Range("R2").Value = 9.338605659243999
$ws.Range("S2").Value = 0.008328107259900691
$ws.Range("T2").Value = 0.008328107259900689

# Row 3
$ws.Range("I3").Value = 0.2503572190582515
$ws.Range("J3").Value = 0.2503572190582515
$ws.Range("M3").Value = 1.626140333333333
$ws.Range("N3").Value = 4.878420999999999
$ws.Range("O3").Value = 0.1053568287437347
$ws.Range("P3").Value = 0.1053568287437347
$ws.Range("Q3").Value = 3.286366736240444
$ws.Range("S3").Value = 0.02637684265307788
$ws.Range("T3").Value = 0.02637684265307788

# Row 4
$ws.Range("I4").Value = 0.2503572190582515
$ws.Range("J4").Value = 0.2503572190582515
$ws.Range("M4").Value = 10.254745
$ws.Range("N4").Value = 30.764235
$ws.Range("O4").Value = 0.6643998618255804
$ws.Range("P4").Value = 0.6643998618255803
$ws.Range("Q4").Value = 20.72444312819334
$ws.Range("R4").Value = 186.51998815374
$ws.Range("S4").Value = 0.1663373017493389
$ws.Range("T4").Value = 0.1663373017493389

# Row 5
$ws.Range("I5").Value = 0.2503572190582515
$ws.Range("J5").Value = 0.2503572190582515
$ws.Range("M5").Value = 3.040282666666667
$ws.Range("N5").Value = 9.120848000000001
$ws.Range("O5").Value = 0.1969784118126819
$ws.Range("P5").Value = 0.1969784118126819
$ws.Range("Q5").Value = 6.14429371173689
$ws.Range("R5").Value = 55.29864340563201
$ws.Range("S5").Value = 0.04931496739593408
$ws.Range("T5").Value = 0.04931496739593408

# Row 6
$ws.Range("G6").Value = 6.051349666666667
$ws.Range("H6").Value = 18.154049
$ws.Range("I6").Value = 0.7496427809417484
$ws.Range("J6").Value = 0.7496427809417485
$ws.Range("M6").Value = 0.5134303333333333
$ws.Range("N6").Value = 1.540291
$ws.Range("O6").Value = 0.03326489761800302
$ws.Range("P6").Value = 0.03326489761800301
$ws.Range("Q6").Value = 3.106946476473222
$ws.Range("R6").Value = 27.962518288259
$ws.Range("S6").Value = 0.02493679035810233
$ws.Range("T6").Value = 0.02493679035810232

# Row 7
$ws.Range("G7").Value = 6.051349666666667
$ws.Range("H7").Value = 18.154049
$ws.Range("I7").Value = 0.7496427809417484
$ws.Range("J7").Value = 0.7496427809417485
$ws.Range("M7").Value = 1.626140333333333
$ws.Range("N7").Value = 4.878420999999999
$ws.Range("O7").Value = 0.1053568287437347
$ws.Range("P7").Value = 0.1053568287437347
$ws.Range("Q7").Value = 9.840343764069887
$ws.Range("R7").Value = 88.563093876629
$ws.Range("S7").Value = 0.07897998609065682
$ws.Range("T7").Value = 0.07897998609065683

# Row 8
$ws.Range("G8").Value = 6.051349666666667
$ws.Range("H8").Value = 18.154049
$ws.Range("I8").Value = 0.7496427809417484
$ws.Range("J8").Value = 0.7496427809417485
$ws.Range("M8").Value = 10.254745
$ws.Range("N8").Value = 30.764235
$ws.Range("O8").Value = 0.6643998618255804
$ws.Range("P8").Value = 0.6643998618255803
$ws.Range("Q8").Value = 62.05504773750167
$ws.Range("R8").Value = 558.495429637515
$ws.Range("S8").Value = 0.4980625600762415
$ws.Range("T8").Value = 0.4980625600762415

# Row 9
$ws.Range("G9").Value = 6.051349666666667
$ws.Range("H9").Value = 18.154049
$ws.Range("I9").Value = 0.7496427809417484
$ws.Range("J9").Value = 0.7496427809417485
$ws.Range("M9").Value = 3.040282666666667
$ws.Range("N9").Value = 9.120848000000001
$ws.Range("O9").Value = 0.1969784118126819
$ws.Range("P9").Value = 0.1969784118126819
$ws.Range("Q9").Value = 18.39781350150578
$ws.Range("R9").Value = 165.580321513552
$ws.Range("S9").Value = 0.1476634444167478
$ws.Range("T9").Value = 0.1476634444167478
